$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Swap-CellText($table, $row, $colA, $colB) {
    $cellA = $table.Cell($row, $colA)
    $cellB = $table.Cell($row, $colB)
    $textA = $cellA.Range.Text
    $textB = $cellB.Range.Text
    $cellA.Range.Text = $textB
    $cellB.Range.Text = $textA
}

# Header row (Yes/No labels under Ecology / Social): swap each No/Yes pair
Swap-CellText $t 2 2 3
Swap-CellText $t 2 4 5

# Data row (counts under Ecology / Social): swap each pair of counts to match
Swap-CellText $t 3 2 3
Swap-CellText $t 3 4 5
